$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder header row: Nombre / Apellido / Correo / Telefono now occupy B1:E1
$ws.Range("B1").Value = "Nombre"
$ws.Range("C1").Value = "Apellido"
$ws.Range("D1").Value = "Correo"
$ws.Range("E1").Value = "Teléfono"

# Update the active selection to D3
$ws.Range("D3").Select()
